$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.093.89"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.970.04"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "328.99"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "1.009"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4988"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "0.4218"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "0.09253"
$ws.Range("E10").Value = "  +5.19%  "
$ws.Range("D11").Value = "1.100"
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("D12").Value = "22.86"
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").Value = "1.980.68"
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").Value = "7.910"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("D15").Value = "6.468"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "1.011"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.00001106"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "91.90"
$ws.Range("E18").Value = "  -4.50%  "
$ws.Range("D19").Value = "0.06741"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").Value = "19.29"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "5.966"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "29.124.44"
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("D24").Value = "11.95"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "2.271"
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("D26").Value = "2.220.06"
$ws.Range("E26").Value = "  -6.46%  "
$ws.Range("D27").Value = "20.71"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").Value = "155.54"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("D29").Value = "6.324"
$ws.Range("E29").Value = "  -5.24%  "
$ws.Range("D30").Value = "2.260"
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").Value = "126.85"
$ws.Range("D32").Value = "1.048"
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "0.09863"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("D34").Value = "1.517"
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("D35").Value = "5.827"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").Value = "3.725"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "0.02436"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").Value = "9.052"
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").Value = "0.06398"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "0.6482"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "11.45"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "0.1996"
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("D44").Value = "1.008"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "0.6221"
$ws.Range("E45").Value = "  -1.42%  "
$ws.Range("D46").Value = "1.358"
$ws.Range("E46").Value = "  +8.77%  "
$ws.Range("D47").Value = "2.202"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "13.29"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "3.477"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").Value = "0.00000000325"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "0.06973"
$ws.Range("E51").Value = "  -0.42%  "
